{"js": "// Load all paragraphs in the document body so we can locate the four\n// \"answer\" paragraphs that the commit filled in.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// There are two identical \"What is the type of the query output in the\n// above question?\" paragraphs (Q4/Q5 answer and Q6/Q7 answer) - grab them\n// in document order.\nconst typeQuestionParas = [];\nlet extensionMethodPara = null;\nlet iGroupingPara = null;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t === \"English: What is the type of the query output in the above question?\") {\n    typeQuestionParas.push(items[i]);\n  } else if (t.indexOf(\"still work to print the list of\") !== -1) {\n    extensionMethodPara = items[i];\n  } else if (t.indexOf(\"IGrouping\") !== -1) {\n    iGroupingPara = items[i];\n  }\n}\n\nif (typeQuestionParas.length < 2) {\n  throw new Error(\"Could not find both 'What is the type...' paragraphs\");\n}\nif (!extensionMethodPara) {\n  throw new Error(\"Could not find the Q1 extension method paragraph\");\n}\nif (!iGroupingPara) {\n  throw new Error(\"Could not find the IGrouping paragraph\");\n}\n\nconst firstTypeQuestion = typeQuestionParas[0];\nconst secondTypeQuestion = typeQuestionParas[1];\n\n// Q4/Q5 answer: \"... question? IOrderedEnumerable.\"\nfirstTypeQuestion.getRange(\"End\").insertText(\" IOrderedEnumerable.\", Word.InsertLocation.replace);\n\n// Q6/Q7 answer: \"... question? IEnumerable\"\nsecondTypeQuestion.getRange(\"End\").insertText(\" IEnumerable\", Word.InsertLocation.replace);\n\n// Q8 (why doesn't the extension method work) answer.\nextensionMethodPara.getRange(\"End\").insertText(\n  \" The typing of the result is different.\",\n  Word.InsertLocation.replace\n);\n\n// IGrouping vs Dictionary answer.\niGroupingPara.getRange(\"End\").insertText(\n  \" They both store objects in collections. The difference is that IGrouping stores objects that have a key in common whereas the dictionary stores objects in key/value pairs.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// The author's cursor (tracked by Word's built-in \"_GoBack\" bookmark) ended\n// up in the middle of the word \"whereas\" - relocate the bookmark there to\n// match (it originally sat right before the paragraph's second sentence).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst whereResults = context.document.body.search(\"common where\", { matchCase: true });\nwhereResults.load(\"items\");\nawait context.sync();\nif (whereResults.items.length > 0) {\n  whereResults.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# There are two identical \"English: What is the type of the query output in\n# the above question?\" paragraphs in the document (one after the\n# alphabetical-order question, one after the FirstOrDefault question).\n# Find.Execute walks forward through the document, so calling it twice in a\n# row (continuing the search range from where the previous hit ended) lands\n# on the first occurrence and then the second occurrence in turn.\n\n# --- 1st occurrence: append \" IOrderedEnumerable.\" ---\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"What is the type of the query output in the above question?\") | Out-Null\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\" IOrderedEnumerable.\")\n$afterFirst = $rng.End\n\n# --- 2nd occurrence: append \" IEnumerable\" ---\n$rng2 = $d.Range($afterFirst, $d.Content.End)\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"What is the type of the query output in the above question?\") | Out-Null\n$rng2.Collapse(0)  # wdCollapseEnd\n$rng2.InsertAfter(\" IEnumerable\")\n\n# --- \"...still work to print the list of students?\" -> append explanation ---\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute(\"still work to print the list of students?\") | Out-Null\n$rng3.Collapse(0)  # wdCollapseEnd\n$rng3.InsertAfter(\" The typing of the result is different.\")\n\n# --- IGrouping<T, K> vs Dictionary<T, K> question -> append full answer ---\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Execute(\"How are they different?\") | Out-Null\n$rng4.Collapse(0)  # wdCollapseEnd\n$rng4.InsertAfter(\" They both store objects in collections. The difference is that IGrouping stores objects that have a key in common whereas the dictionary stores objects in key/value pairs.\")\n\n# The author's cursor (tracked by Word's built-in \"_GoBack\" bookmark) ended up\n# in the middle of the word \"whereas\" when they saved - relocate the bookmark\n# there to match (it originally sat right before the paragraph's 2nd sentence).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$rngB = $d.Content\n$rngB.Find.ClearFormatting()\n$rngB.Find.Execute(\"common where\") | Out-Null\n$rngB.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rngB) | Out-Null\n"}
